# Append the new "Realtime Sync / OT-like Merge / Live Cursor Update" status
# section to the end of the document: a blank spacer, a "---" separator, the
# section heading, an "Updated:" line, another blank spacer, the column
# header row, and three new module rows (Realtime Sync (WebSocket), OT/CRDT
# Conflict Handling, Live Cursor Co-editing). Each multi-column row is a
# single paragraph whose cells are separated by literal tab runs, matching
# the style already used by the "Cumulative Main Requirement Status Matrix"
# section earlier in this document.

$d = $word.ActiveDocument

# The run formatting (font + size) shared by every paragraph in this
# document, reused verbatim for the new paragraphs so they match style.
function New-RunProps() {
    return '<w:rPr><w:rFonts w:ascii="Helvetica Light" w:hAnsi="Helvetica Light" w:cs="Helvetica Light"/><w:sz w:val="24"/><w:sz-cs w:val="24"/></w:rPr>'
}

# Builds the <w:t>/<w:tab/> run contents for one paragraph from an array of
# cell strings (a 1-element array yields plain text with no tabs).
function New-CellsXml($Cells) {
    $out = ''
    for ($i = 0; $i -lt $Cells.Count; $i++) {
        if ($i -gt 0) { $out += '<w:tab/>' }
        $out += '<w:t xml:space="preserve">' + $Cells[$i] + '</w:t>'
    }
    return $out
}

# Builds a full <w:p> paragraph (empty pPr + single run) for one row.
function New-ParagraphXml($Cells) {
    $rPr = New-RunProps
    $cellsXml = New-CellsXml $Cells
    return '<w:p><w:pPr/><w:r>' + $rPr + $cellsXml + '</w:r></w:p>'
}

$newParagraphs = @(
    (New-ParagraphXml @('')),
    (New-ParagraphXml @('---')),
    (New-ParagraphXml @('Realtime Sync / OT-like Merge / Live Cursor Update')),
    (New-ParagraphXml @('Updated: 2026-02-18')),
    (New-ParagraphXml @('')),
    (New-ParagraphXml @('Module Name', 'Developed', 'Partial Developed', 'Need To Develop')),
    (New-ParagraphXml @(
        'Realtime Sync (WebSocket)',
        'Added websocket collaboration endpoint `/ws/storefront/{storeId}` with snapshot/presence/op messaging and client heartbeat integration in Store Builder',
        'In-memory room state only; no distributed backplane yet',
        'Redis/backplane fanout for multi-instance horizontal scale'
    )),
    (New-ParagraphXml @(
        'OT/CRDT Conflict Handling',
        'Implemented revision-based operation acceptance + conflict response (`conflict` returns server snapshot) for safe convergence',
        'This is OT-like revision gating, not full CRDT/OT transform algorithm',
        'True CRDT/OT transform engine with operation-level semantic merge'
    )),
    (New-ParagraphXml @(
        'Live Cursor Co-editing',
        'Added live cursor broadcasts (`cursor` events) and remote cursor presence indicators in builder UI',
        'Cursor coordinates currently minimal (node-focused, not pixel cursor map)',
        'Full live cursor coordinates, colored user avatars, viewport-aware cursor rendering'
    ))
)

$newSectionXml = $newParagraphs -join ''

# The fragment needs the WordprocessingML namespace declared somewhere on
# it; put it on the first element only (namespace applies to descendants).
$newSectionXml = $newSectionXml -replace '^<w:p>', '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">'

# A zero-length Range collapsed to the very end of the main document story:
# inserting OOXML into it appends after the last existing paragraph (rather
# than overwriting/replacing it) and lands the new content right before
# </w:body> (ahead of the sectPr), exactly where the diff adds it.
$endRange = $d.Content
$endRange.Collapse(0)

[void]$endRange.InsertXML($newSectionXml)
